# Machine Learning on Original Imbalanced Dataset
# Adds a new "Sheet5" with model-comparison metrics (Accuracy / Precision /
# Recall / F1 / ROC AUC) for five classifiers, appended after Sheet4, and
# makes it the active/selected sheet (matching the target workbook state).

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the LAST tab (after Sheet4) -------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws5 = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws5.Name = "Sheet5"

# --- Header row (row 2) ------------------------------------------------------
$ws5.Range("B2").Value = "Model"
$ws5.Range("C2").Value = "Accuracy"
$ws5.Range("D2").Value = "Precision"
$ws5.Range("E2").Value = "Recall"
$ws5.Range("F2").Value = "F1 Score"
$ws5.Range("G2").Value = "ROC AUC Score"

# --- Data rows (rows 3-7) ----------------------------------------------------
$ws5.Range("B3").Value = "Logistic Regression"
$ws5.Range("C3").Value = 0.78
$ws5.Range("D3").Value = 0.79
$ws5.Range("E3").Value = 0.97
$ws5.Range("F3").Value = 0.87
$ws5.Range("G3").Value = 0.55

$ws5.Range("B4").Value = "Decision Tree"
$ws5.Range("C4").Value = 0.76
$ws5.Range("D4").Value = 0.85
$ws5.Range("E4").Value = 0.85
$ws5.Range("F4").Value = 0.85
$ws5.Range("G4").Value = 0.66

$ws5.Range("B5").Value = "Random Forest"
$ws5.Range("C5").Value = 0.78
$ws5.Range("D5").Value = 0.84
$ws5.Range("E5").Value = 0.88
$ws5.Range("F5").Value = 0.86
$ws5.Range("G5").Value = 0.65

$ws5.Range("B6").Value = "Gradient Boosting"
$ws5.Range("C6").Value = 0.8
$ws5.Range("D6").Value = 0.82
$ws5.Range("E6").Value = 0.96
$ws5.Range("F6").Value = 0.88
$ws5.Range("G6").Value = 0.6

$ws5.Range("B7").Value = "LightGBM"
$ws5.Range("C7").Value = 0.81
$ws5.Range("D7").Value = 0.84
$ws5.Range("E7").Value = 0.93
$ws5.Range("F7").Value = 0.88
$ws5.Range("G7").Value = 0.65

# --- Formatting: reuse the same header/body styles already used on Sheet4 ---
# (bold + thin box border for the header row, thin box border for the body)
$headerRange = $ws5.Range("B2:G2")
$headerRange.Borders.Color = 0
$headerRange.Borders.LineStyle = 1
$headerRange.Font.Bold = $true

$bodyRange = $ws5.Range("B3:G7")
$bodyRange.Borders.Color = 0
$bodyRange.Borders.LineStyle = 1

# --- Column widths (best-fit-like, matches target sheet look) ---------------
$ws5.Columns.Item(2).ColumnWidth = 16.61
$ws5.Columns.Item(3).ColumnWidth = 7.83
$ws5.Columns.Item(5).ColumnWidth = 5.28
$ws5.Columns.Item(6).ColumnWidth = 7.17
$ws5.Columns.Item(7).ColumnWidth = 13.17

# --- Selection / view state matching the target file -------------------------
$ws5.Range("B2:G7").Select()

# Make the new sheet the active tab (moves tabSelected/activeTab here and off
# of Sheet4, matching the diff).
$ws5.Activate()
